# Add the 2022-Q3 sheet (fund-holdings detail) right after "总计",
# shifting the existing quarter sheets along, and insert the matching
# summary row at the top of the "总计" (total) sheet.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)
$q2Sheet = $wb.Worksheets.Item(2)

# --- 1. Duplicate the existing "2022-Q2" sheet so the new sheet starts
#     with identical headers / column formatting, then rename + reposition it.
#     Worksheet.Copy(Before) drops the clone immediately before $q2Sheet, i.e.
#     right after the "总计" sheet - exactly where "2022-Q3" belongs.
$q2Sheet.Copy($q2Sheet)
$q3Sheet = $wb.Worksheets.Item(2)
$q3Sheet.Name = "2022-Q3"

# --- 2. Extend the cloned sheet with the 10 fund rows for 2022-Q3 (template
#     only carried 4 rows). Column A keeps the bold/boxed index style used by
#     the other cells in column A, so clone that style down first.
$q3Sheet.Range("A5").Copy($q3Sheet.Range("A6:A11"))

# --- 3. Columns B, D, E, F, G hold text that looks numeric ("004616", "55.17",
#     ...) and must stay text (leading zeros, exact decimal strings) instead of
#     being auto-converted to numbers, so force Text format before writing.
$q3Sheet.Range("B2:B11").NumberFormat = "@"
$q3Sheet.Range("D2:G11").NumberFormat = "@"

$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "501208"
$q3Sheet.Range("C2").Value = "中欧创新未来混合（LOF）"
$q3Sheet.Range("D2").Value = "55.17"
$q3Sheet.Range("E2").Value = "85.30"
$q3Sheet.Range("F2").Value = "3.04"
$q3Sheet.Range("G2").Value = "1.6772"
$q3Sheet.Range("H2").Value = 9

$q3Sheet.Range("A3").Value = 1
$q3Sheet.Range("B3").Value = "004616"
$q3Sheet.Range("C3").Value = "中欧电子信息产业沪港深股票A"
$q3Sheet.Range("D3").Value = "5.01"
$q3Sheet.Range("E3").Value = "92.97"
$q3Sheet.Range("F3").Value = "5.80"
$q3Sheet.Range("G3").Value = "0.2906"
$q3Sheet.Range("H3").Value = 3

$q3Sheet.Range("A4").Value = 2
$q3Sheet.Range("B4").Value = "005763"
$q3Sheet.Range("C4").Value = "中欧电子信息产业沪港深股票C"
$q3Sheet.Range("D4").Value = "3.88"
$q3Sheet.Range("E4").Value = "92.97"
$q3Sheet.Range("F4").Value = "5.80"
$q3Sheet.Range("G4").Value = "0.2250"
$q3Sheet.Range("H4").Value = 3

$q3Sheet.Range("A5").Value = 3
$q3Sheet.Range("B5").Value = "001411"
$q3Sheet.Range("C5").Value = "诺安创新驱动灵活配置混合A"
$q3Sheet.Range("D5").Value = "3.98"
$q3Sheet.Range("E5").Value = "80.56"
$q3Sheet.Range("F5").Value = "3.55"
$q3Sheet.Range("G5").Value = "0.1413"
$q3Sheet.Range("H5").Value = 7

$q3Sheet.Range("A6").Value = 4
$q3Sheet.Range("B6").Value = "002051"
$q3Sheet.Range("C6").Value = "诺安创新驱动灵活配置混合C"
$q3Sheet.Range("D6").Value = "1.80"
$q3Sheet.Range("E6").Value = "80.56"
$q3Sheet.Range("F6").Value = "3.55"
$q3Sheet.Range("G6").Value = "0.0639"
$q3Sheet.Range("H6").Value = 7

$q3Sheet.Range("A7").Value = 5
$q3Sheet.Range("B7").Value = "001097"
$q3Sheet.Range("C7").Value = "华泰柏瑞积极优选股票A"
$q3Sheet.Range("D7").Value = "1.20"
$q3Sheet.Range("E7").Value = "83.01"
$q3Sheet.Range("F7").Value = "3.05"
$q3Sheet.Range("G7").Value = "0.0366"
$q3Sheet.Range("H7").Value = 5

$q3Sheet.Range("A8").Value = 6
$q3Sheet.Range("B8").Value = "562520"
$q3Sheet.Range("C8").Value = "华夏中证智选1000成长创新策略ETF"
$q3Sheet.Range("D8").Value = "0.44"
$q3Sheet.Range("E8").Value = "96.91"
$q3Sheet.Range("F8").Value = "1.19"
$q3Sheet.Range("G8").Value = "0.0052"
$q3Sheet.Range("H8").Value = 2

$q3Sheet.Range("A9").Value = 7
$q3Sheet.Range("B9").Value = "016283"
$q3Sheet.Range("C9").Value = "华泰柏瑞积极优选股票C"
$q3Sheet.Range("D9").Value = "0.16"
$q3Sheet.Range("E9").Value = "83.01"
$q3Sheet.Range("F9").Value = "3.05"
$q3Sheet.Range("G9").Value = "0.0049"
$q3Sheet.Range("H9").Value = 5

$q3Sheet.Range("A10").Value = 8
$q3Sheet.Range("B10").Value = "167703"
$q3Sheet.Range("C10").Value = "德邦量化优选股票（LOF）C"
$q3Sheet.Range("D10").Value = "0.30"
$q3Sheet.Range("E10").Value = "90.17"
$q3Sheet.Range("F10").Value = "1.03"
$q3Sheet.Range("G10").Value = "0.0031"
$q3Sheet.Range("H10").Value = 3

$q3Sheet.Range("A11").Value = 9
$q3Sheet.Range("B11").Value = "167702"
$q3Sheet.Range("C11").Value = "德邦量化优选股票（LOF）A"
$q3Sheet.Range("D11").Value = "0.17"
$q3Sheet.Range("E11").Value = "90.17"
$q3Sheet.Range("F11").Value = "1.03"
$q3Sheet.Range("G11").Value = "0.0018"
$q3Sheet.Range("H11").Value = 3

# --- 4. Drop the leftover "Text" number-format flag the writes above applied
#     (copying the format of a pristine, never-touched cell resets it to the
#     plain/default style the other text cells already use).
$blank = $q3Sheet.Range("Z100")
$blank.Copy()
$q3Sheet.Range("B2:B11").PasteSpecial(-4122)
$q3Sheet.Range("D2:G11").PasteSpecial(-4122)
$q3Sheet.Application.CutCopyMode = $false

# --- 5. Insert the matching 2022-Q3 summary row at the top of "总计",
#     pushing the older quarters down by one row (same index-column style
#     trick as above).
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A6").Copy($totalSheet.Range("A2"))
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 10
$totalSheet.Range("D2").Value = 2.45

# --- 6. Renumber the A-column index (0-based) for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

# --- 7. Leave the view the way it started, with "总计" selected.
$totalSheet.Activate()
$wb.Application.CutCopyMode = $false
